$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.372.46"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").Value = "1.667.22"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'218.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").Value = "'1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("D9").Value = "'0.2567"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.79%  "
$ws.Range("D10").Value = "'19.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").Value = "'0.07663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.678.87"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D14").Value = "1.896.52"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").Value = "0.0₅8042"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "'64.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.98%  "
$ws.Range("D18").Value = "26.392.84"
$ws.Range("E18").Value = "  -3.57%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'211.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "'4.406"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.79%  "
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("D23").Value = "'5.903"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'144.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").Value = "'0.1162"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").Value = "'6.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("D30").Value = "'0.05259"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").Value = "'3.379"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D33").Value = "'3.220"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.22%  "
$ws.Range("D34").Value = "'1.567"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("D36").Value = "'2.377"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("D37").Value = "'0.9286"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "1.154.23"
$ws.Range("E39").Value = "  +10.56%  "
$ws.Range("D40").Value = "'0.01597"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").Value = "'0.8480"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "'1.007"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "'99.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "1.806.03"
$ws.Range("E46").Value = "  -5.40%  "
$ws.Range("D47").Value = "'0.4499"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'56.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'7.939"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "'0.05101"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
